$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row 45 - MindSensorPressureSensor: claim the sensor (SampleProvider interface, Pressure mode)
$ws.Range("D45").Value = "Lawrie"
$ws.Range("E45").Value = "N"
$ws.Range("F45").Value = "Pressure"
$ws.Range("G45").Value = "SampleProvider"

# Row 69 - RCXTemperatureSensor: claim the sensor (SampleProvider interface, Temperature mode)
$ws.Range("D69").Value = "Lawrie"
$ws.Range("E69").Value = "N"
$ws.Range("F69").Value = "Temperature"
$ws.Range("G69").Value = "SampleProvider"

# Move the active selection in the frozen bottom-right pane to I66
$null = $ws.Range("I66").Select()
